# "Small mistake on one of the tasks" - fix a typo in the Sprint 4 task list:
# "Tiago, Afonso, and João ..." should read "Tiago, Pedro, and João ..."
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Folha1")

$ws.Range("B4").Value = "Tiago, Pedro, and João should analyze the code and identify classes and data Structures important to implement trash."

# Leave the edited cell selected/active, matching the author's final cursor position.
$ws.Range("B4").Select()
